# Update "想去人数" (want-to-go count) values in column F
# on the "展览" and "全部类型" sheets to reflect refreshed scrape data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8875
    3  = 8314
    7  = 48
    8  = 150
    9  = 160
    13 = 209
    14 = 5391
    17 = 87
    18 = 17
    19 = 24
    21 = 159
    22 = 189
    23 = 13
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
